$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log rows appended to "Registro de Actividad" (rows 107-120).
# Each entry: Timestamp (A), Accion (D), Detalles Adicionales (E)
$entries = @(
    @("2025-06-06 19:47:06", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:47:13", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:47:46", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:47:57", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:49:07", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:49:18", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:49:35", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:49:52", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:50:14", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 19:53:19", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 20:02:25", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 20:02:26", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 20:04:52", "Cierre Aplicacion", "Cerrado desde login."),
    @("2025-06-06 20:08:23", "Cierre Aplicacion", "Cerrado desde login.")
)

$startRow = 107
for ($i = 0; $i -lt $entries.Count; $i++) {
    $row = $startRow + $i
    $entry = $entries[$i]
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 4).Value = $entry[1]
    $ws.Cells.Item($row, 5).Value = $entry[2]
}
